$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1177.5
$ws.Range("I92").Value = 434.625
$ws.Range("J92").Value = 2663.25
$ws.Range("K92").Value = 434.625
$ws.Range("L92").Value = 2663.25
$ws.Range("M92").Value = 813.375
$ws.Range("N92").Value = -5159.25

$ws.Range("H111").Value = 23812934
$ws.Range("I111").Value = 30304916
$ws.Range("K111").Value = 90914748
$ws.Range("M111").Value = -90911681

$ws.Range("H138").Value = 313384.22
$ws.Range("I138").Value = 3497.625
$ws.Range("J138").Value = 384215.44
$ws.Range("K138").Value = 10492.875
$ws.Range("L138").Value = 1152646.32
$ws.Range("M138").Value = -5352.875
$ws.Range("N138").Value = -1162926.32

$ws.Range("H141").Value = 2166169.2
$ws.Range("I141").Value = 1175.2609
$ws.Range("K141").Value = 3525.7827
$ws.Range("M141").Value = 1654.2173

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2706133.8
$ws.Range("I32").Value = 4504.4287
$ws.Range("K32").Value = 4504.4287
$ws.Range("M32").Value = -4217.4287

$ws.Range("H61").Value = 1429.3939
$ws.Range("I61").Value = 1446.25
$ws.Range("J61").Value = 890
$ws.Range("K61").Value = 1446.25
$ws.Range("L61").Value = 890
$ws.Range("M61").Value = -1234.25
$ws.Range("N61").Value = -1314

$ws.Range("H74").Value = 1020.1429
$ws.Range("I74").Value = 1009.2222
$ws.Range("J74").Value = 1085.6666
$ws.Range("K74").Value = 1009.2222
$ws.Range("L74").Value = 1085.6666
$ws.Range("M74").Value = -135.2222
$ws.Range("N74").Value = -2833.6666

$ws.Range("H77").Value = 1020.1429
$ws.Range("I77").Value = 1009.2222
$ws.Range("J77").Value = 1085.6666
$ws.Range("K77").Value = 5046.111
$ws.Range("L77").Value = 5428.333000000001
$ws.Range("M77").Value = -678.1109999999999
$ws.Range("N77").Value = -14164.333

$ws.Range("H136").Value = 1429.3939
$ws.Range("I136").Value = 1446.25
$ws.Range("J136").Value = 890
$ws.Range("K136").Value = 4338.75
$ws.Range("L136").Value = 2670
$ws.Range("M136").Value = -1788.75
$ws.Range("N136").Value = -7770

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 12349505
$ws.Range("I20").Value = 23815562
$ws.Range("J20").Value = 1445.1538
$ws.Range("K20").Value = 23815562
$ws.Range("L20").Value = 1445.1538
$ws.Range("M20").Value = -23815315
$ws.Range("N20").Value = -1939.1538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2027.925
$ws.Range("I31").Value = 1597.6154
$ws.Range("J31").Value = 2827.0715
$ws.Range("K31").Value = 1597.6154
$ws.Range("L31").Value = 2827.0715
$ws.Range("M31").Value = -1302.6154
$ws.Range("N31").Value = -3417.0715

$ws.Range("H34").Value = 2027.925
$ws.Range("I34").Value = 1597.6154
$ws.Range("J34").Value = 2827.0715
$ws.Range("K34").Value = 1597.6154
$ws.Range("L34").Value = 2827.0715
$ws.Range("M34").Value = -1395.6154
$ws.Range("N34").Value = -3231.0715

$ws.Range("H58").Value = 1284.1936
$ws.Range("I58").Value = 1275.28
$ws.Range("J58").Value = 1321.3334
$ws.Range("K58").Value = 1275.28
$ws.Range("L58").Value = 1321.3334
$ws.Range("M58").Value = -1072.28
$ws.Range("N58").Value = -1727.3334

$ws.Range("H99").Value = 1741.6
$ws.Range("I99").Value = 1616.5714
$ws.Range("J99").Value = 2033.3334
$ws.Range("K99").Value = 1616.5714
$ws.Range("L99").Value = 2033.3334
$ws.Range("M99").Value = -118.5714
$ws.Range("N99").Value = -5029.3334

$ws.Range("H126").Value = 1741.6
$ws.Range("I126").Value = 1616.5714
$ws.Range("J126").Value = 2033.3334
$ws.Range("K126").Value = 4849.7142
$ws.Range("L126").Value = 6100.0002
$ws.Range("M126").Value = -2379.7142
$ws.Range("N126").Value = -11040.0002

$ws.Range("H127").Value = 53328
$ws.Range("J127").Value = 53328
$ws.Range("L127").Value = 53328
$ws.Range("N127").Value = -63248

$ws.Range("H136").Value = 1284.1936
$ws.Range("I136").Value = 1275.28
$ws.Range("J136").Value = 1321.3334
$ws.Range("K136").Value = 3825.84
$ws.Range("L136").Value = 3964.0002
$ws.Range("M136").Value = -1275.84
$ws.Range("N136").Value = -9064.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 603.25714
$ws.Range("I107").Value = 328.83334
$ws.Range("J107").Value = 1202
$ws.Range("K107").Value = 986.5000200000001
$ws.Range("L107").Value = 3606
$ws.Range("M107").Value = 933.4999799999999
$ws.Range("N107").Value = -7446

$ws.Range("H131").Value = 900.09
$ws.Range("I131").Value = 425.7143
$ws.Range("J131").Value = 935.7957
$ws.Range("K131").Value = 1277.1429
$ws.Range("L131").Value = 2807.3871
$ws.Range("M131").Value = 3762.8571
$ws.Range("N131").Value = -12887.3871

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 25100.49
$ws.Range("I70").Value = 28345.455
$ws.Range("J70").Value = 4703.5713
$ws.Range("K70").Value = 28345.455
$ws.Range("L70").Value = 4703.5713
$ws.Range("M70").Value = -28075.455
$ws.Range("N70").Value = -5243.5713

$ws.Range("H73").Value = 25100.49
$ws.Range("I73").Value = 28345.455
$ws.Range("J73").Value = 4703.5713
$ws.Range("K73").Value = 28345.455
$ws.Range("L73").Value = 4703.5713
$ws.Range("M73").Value = -27409.455
$ws.Range("N73").Value = -6575.5713

$ws.Range("H113").Value = 1269.8214
$ws.Range("I113").Value = 1120.5883
$ws.Range("J113").Value = 1500.4546
$ws.Range("K113").Value = 1120.5883
$ws.Range("L113").Value = 1500.4546
$ws.Range("M113").Value = 1049.4117
$ws.Range("N113").Value = -5840.4546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10002010
$ws.Range("I7").Value = 12501637
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 12501637
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -12501525
$ws.Range("N7").Value = -3724

$ws.Range("H122").Value = 2850.3572
$ws.Range("I122").Value = 1882.2222
$ws.Range("J122").Value = 4593
$ws.Range("K122").Value = 5646.6666
$ws.Range("L122").Value = 13779
$ws.Range("M122").Value = -3196.6666
$ws.Range("N122").Value = -18679

$ws.Range("H126").Value = 10002010
$ws.Range("I126").Value = 12501637
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 37504911
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -37502441
$ws.Range("N126").Value = -15440
